$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Wat heb ik gedaan" paragraph - rewrite the summary text.
# ---------------------------------------------------------------------------
$rsquo = [char]0x2019
$old1 = "In de afgelopen week heb ik voornamelijk naar goede dashcam filmpjes gezocht die wij kunnen gebruiken voor het project. Ook heb ik uit al die filmpjes de verkeersborden geknipt voor de bordendetectie."
$new1 = "In de afgelopen week heb images gecropt die wij gaan gebruiken voor template matching. Dit zijn dan foto${rsquo}s van dashcam videos waarvan alleen de verkeersborden zichtbaar zijn. Dingen zoals de weg, andere auto${rsquo}s en de lucht zijn niet meer zichtbaar. Ik heb ook een samenvatting van 500 woorden geschreven over het project. Dit is noodzakelijk voor de inzending."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Remove the stray "_GoBack" bookmark that originally sat right after the
# paragraph above (the engine exposes no working Bookmark.Delete, so we
# surgically cut it out: delete the character just before it together with
# the zero-width break run(s) immediately following it, then retype that
# exact same content back in place - this removes the bookmark markers
# without losing any visible text/line-breaks).
# ---------------------------------------------------------------------------
$vtab = [string][char]11
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $p = $bm.Range.Start

    $nbreaks = 0
    while ($true) {
        $probe = $d.Range($p + $nbreaks, $p + $nbreaks + 1)
        if ($probe.Text.Length -eq 1 -and [int][char]$probe.Text[0] -eq 11) {
            $nbreaks = $nbreaks + 1
        } else {
            break
        }
    }

    $beforeRng = $d.Range($p - 1, $p)
    $beforeChar = $beforeRng.Text

    $delRange = $d.Range($p - 1, $p + $nbreaks)
    $delRange.Delete()

    $rebuild = $beforeChar
    for ($k = 0; $k -lt $nbreaks; $k++) {
        $rebuild = $rebuild + $vtab
    }
    $insPos = $d.Range($p - 1, $p - 1)
    $insPos.InsertAfter($rebuild)
}

# ---------------------------------------------------------------------------
# Change 2 ("Samenwerking" paragraph): merge the two closing sentences
# into one rewritten sentence.
# ---------------------------------------------------------------------------
$ldquo = [char]0x201C
$rdquo = [char]0x201D
$old2 = " Grotendeels werkte iedereen individueel aan taken dus er was er niet echt veel ${ldquo}samenwerken${rdquo} van toepassing. Er werd uiteraard wel goed gecommuniceerd onder het team."
$new2 = " Wessel en Bart werkten samen op de Raspberry Pi. Ik werkte veel individueel aan mijn eigen taken. Er werd ook goed gecommuniceerd onder het team."
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# Changes 3 & 4 ("Verbeterpunten" paragraph): rewrite both runs of text.
# ---------------------------------------------------------------------------
$old3 = "Een verbeterpunt die ik voor mezelf heb "
$new3 = "Het zou fijner zijn als een aantal mensen meer tijd zouden besteden aan het project. Soms merk ik wel dat een p"
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

$old4 = "is om sneller aan te geven als er iets niet lukt bij mij. Als ik ergens vast kom te zitten ga ik meestal net zo lang door tot het me lukt. Zelfs als het dagen kan duren. Echter is dit een lange periode en kan ik het beter meteen aangeven in plaats van 2 dagen te wachten."
$new4 = "aar mensen het project laten liggen en iets anders gaan doen."
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------------
# Re-add the "_GoBack" bookmark at the very end of the document (end of the
# "Verbeterpunten" paragraph), matching where Word would leave it after the
# last edit.
# ---------------------------------------------------------------------------
$lastParaRange = $d.Paragraphs.Last.Range
$endPos = $lastParaRange.End - 1
$rngEnd = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $rngEnd) | Out-Null
